$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the values in row 2: Testcase ID and Test data/Input data changed,
# Expected values stays the same text ("landed on home page").
# Write order matters for shared-string table index assignment.
$ws.Range("C2").Value = "landed on home page"
$ws.Range("B2").Value = "N/A (No specific input data required)"
$ws.Range("A2").Value = "TC_cura_validatehomepage_004"

# Column width changes (target widths are 33.44140625 / 38.109375 chars of
# stored XML width; the ColumnWidth setter here only supports 1/6-character
# granularity, so we choose the closest reachable value: 98/3 -> 33.5 stored,
# 112/3 -> 38.1667 stored, the nearest quantized results to the target).
$ws.Columns.Item(1).ColumnWidth = 32.6666666666667
$ws.Columns.Item(2).ColumnWidth = 37.3333333333333

# Selection / view changes: scroll so column B is the left-most visible
# column, then select A1:D1 (matches the target sheetView/selection).
$win = $excel.ActiveWindow
$win.ScrollColumn = 2
$win.ScrollRow = 1
$ws.Range("A1:D1").Select()
